$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) and the workbook title to reflect the new "through" date.
$ws.Name = "Through 2021-12-29"

# Update the "December (through 12-xx)" label.
$ws.Range("A14").Value = "December (through 12-29)"

# Row 14 - December (through 12-29) values.
$ws.Range("C14").Value = 40
$ws.Range("D14").Value = 0.09089999999999999
$ws.Range("F14").Value = 85
$ws.Range("G14").Value = 0.08599999999999999
$ws.Range("H14").Value = 13
$ws.Range("I14").Value = 98
$ws.Range("J14").Value = 0.1171
$ws.Range("L14").Value = 67
$ws.Range("M14").Value = 0.0694
$ws.Range("O14").Value = 57
$ws.Range("P14").Value = 0.0806
$ws.Range("R14").Value = 127
$ws.Range("S14").Value = 0.06619999999999999
$ws.Range("U14").Value = 180
$ws.Range("V14").Value = 0.011

# Row 15 - Total values.
$ws.Range("C15").Value = 298
$ws.Range("D15").Value = 0.1104
$ws.Range("F15").Value = 589
$ws.Range("G15").Value = 0.1021
$ws.Range("H15").Value = 76
$ws.Range("I15").Value = 856
$ws.Range("J15").Value = 0.0815
$ws.Range("L15").Value = 675
$ws.Range("M15").Value = 0.1048
$ws.Range("O15").Value = 537
$ws.Range("P15").Value = 0.099
$ws.Range("R15").Value = 1327
$ws.Range("S15").Value = 0.0521
$ws.Range("U15").Value = 1722
$ws.Range("V15").Value = 0.0564
